$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACIONA_DISCADOR_OUTBOUND")

# Insert a new row above the current row 2 (DATA), shifting everything down.
# Copy formatting from the last data row (DATA_INSERT, row 20) so the new
# row inherits the same visual style used for other "merged/IDENTITY" rows.
$ws.Rows.Item(20).Copy()
$ws.Rows.Item(2).Insert()

# Populate the new primary-key row: EVENTO_ID.
$ws.Range("A2").Value = "EVENTO_ID"
$ws.Range("B2").Value = "ID DO ACIONAMENTO NO DW"
$ws.Range("C2").Value = "BIGINT"
$ws.Range("D2").Value = 19
$ws.Range("E2").Value = "PK"
$ws.Range("F2").Value = "IDENTITY"
$ws.Range("G2").Value = $null
$ws.Range("H2").Value = $null
$ws.Range("I2").Value = $null
$ws.Range("J2").Value = $null

# The merged F20:J20 region from the copied row survived the insert as
# F2:J2 - make sure it stays merged.
if (-not $ws.Range("F2:J2").MergeCells) {
    $ws.Range("F2:J2").Merge()
}

# CALL_ID used to be the only primary key (row 7); now that EVENTO_ID is the
# PK, clear the old PK marker. CALL_ID is now row 8 after the insert.
$ws.Range("E8").Value = $null

# Re-apply best-fit widths on the columns whose contents changed.
$ws.Columns.Item("D").AutoFit()
$ws.Columns.Item("E").AutoFit()

# Match the author's final selection.
$ws.Range("E8").Select()
